$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.956.08"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "1.553.91"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.56%  "
$ws.Range("D5").Value = "'206.69"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  +0.49%  "
$ws.Range("D8").Value = "'21.63"
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").Value = "1.775.88"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "1.554.93"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").Value = "26.950.15"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "'61.78"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "'214.33"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("E23").Value = "  +1.41%  "
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("D25").Value = "'153.24"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  +2.05%  "
$ws.Range("D27").Value = "'14.90"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("E29").Value = "  +1.20%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -1.18%  "
$ws.Range("E32").Value = "  +1.43%  "
$ws.Range("D33").Value = "1.382.03"
$ws.Range("E33").Value = "  +1.91%  "
$ws.Range("E34").Value = "  +2.53%  "
$ws.Range("E35").Value = "  +3.21%  "
$ws.Range("D36").Value = "'0.969"
$ws.Range("E36").Value = "  +5.29%  "
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("E38").Value = "  +1.51%  "
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "'0.808"
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("D42").Value = "'0.991"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  +3.02%  "
$ws.Range("D44").Value = "'5.46"
$ws.Range("E44").Value = "  -2.58%  "
$ws.Range("D45").Value = "'63.87"
$ws.Range("E45").Value = "  +1.25%  "
$ws.Range("E46").Value = "  -2.57%  "
$ws.Range("D47").Value = "1.689.42"
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("D48").Value = "'85.91"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").Value = "0.0₇0961"
$ws.Range("E51").Value = "  -1.23%  "
